$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3 (shifts existing rows 3-26 down to 4-27).
# Copy formatting from the row that will end up below (old row 3, now row 4)
# for columns A-C and H, matching Excel's natural "insert row" behavior.
$ws.Rows("3:3").Insert(-4121, 1)

$ws.Range("A4:C4").Copy()
$ws.Range("A3:C3").PasteSpecial(-4122)
$ws.Range("H4").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new task row.
$ws.Range("A3").Value = "Add Stock"
$ws.Range("B3").Value = "Require compartment to be selected"
$ws.Range("C3").Value = 4
$ws.Range("H3").Value = "Guv"

# Bump the "Initial Estimate" for every other existing task by 1.
for ($r = 4; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2 + 1
}

# Update the selection to match the saved view state.
$ws.Range("D14").Select()

# The chart object anchored near the table did not auto-shift with the row
# insert, so nudge it down by one row's height to re-anchor it one row lower.
$co = $ws.ChartObjects(1)
$co.Top = $co.Top + $ws.Rows(3).Height

# Update the chart series range to track the shifted totals row.
$chart = $co.Chart
$chart.SeriesCollection(1).Values = "Sheet1!`$C`$27:`$G`$27"

$excel.CalculateFullRebuild()
